$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dates 2021-05-10 .. 2021-05-13) appended after existing
# row 251 ("aggiornamento fino a 13/03" update).
$newRows = @(
    @{ Row = 252; A = 44326; B = 1;  C = 50; D = 278.1176994103905 },
    @{ Row = 253; A = 44327; B = 11; C = 60; D = 333.7412392924685 },
    @{ Row = 254; A = 44328; B = 2;  C = 61; D = 339.3035932806764 },
    @{ Row = 255; A = 44329; B = 11; C = 55; D = 305.9294693514295 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $prevRow = $row - 1

    # Copy formatting from the row above (column A carries the date style,
    # B/C/D use the default/general style) before writing the new values.
    $ws.Range("A$prevRow").Copy($ws.Range("A$row"))
    $ws.Range("B$prevRow").Copy($ws.Range("B$row"))
    $ws.Range("C$prevRow").Copy($ws.Range("C$row"))
    $ws.Range("D$prevRow").Copy($ws.Range("D$row"))

    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
}
